$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the combustion-chamber temperature input on the "Design" sheet
# (cell B9). Every other changed cell in the workbook (Design!B7, B8,
# B12, B16...B32, plus the dependent formulas on "Printable Results"
# and "Fabrication Stuff") is a formula that recalculates automatically
# from this single input.
$ws.Range("B9").Value = 5642.33

# Reflect the saved cursor position: selection moved to B11 on the
# Design sheet (and the sheet is scrolled back so there's no
# topLeftCell override).
$ws.Range("B11").Select()
